# Update the answers in the "three-digit number divided by one-digit number"
# worksheet to match the newly generated set of problems/answers.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "908÷9=100, 8"; New = "521÷7=74, 3" },
    @{ Old = "618÷9=68, 6";  New = "498÷8=62, 2" },
    @{ Old = "443÷4=110, 3"; New = "268÷4=67, 0" },
    @{ Old = "607÷5=121, 2"; New = "597÷3=199, 0" },
    @{ Old = "995÷5=199, 0"; New = "923÷3=307, 2" },
    @{ Old = "542÷9=60, 2";  New = "319÷5=63, 4" },
    @{ Old = "324÷5=64, 4";  New = "119÷5=23, 4" },
    @{ Old = "866÷5=173, 1"; New = "312÷7=44, 4" },
    @{ Old = "167÷7=23, 6";  New = "624÷8=78, 0" },
    @{ Old = "400÷8=50, 0";  New = "742÷6=123, 4" },
    @{ Old = "649÷7=92, 5";  New = "273÷9=30, 3" },
    @{ Old = "301÷8=37, 5";  New = "923÷4=230, 3" },
    @{ Old = "185÷8=23, 1";  New = "554÷9=61, 5" },
    @{ Old = "430÷6=71, 4";  New = "361÷2=180, 1" },
    @{ Old = "553÷6=92, 1";  New = "457÷9=50, 7" },
    @{ Old = "306÷6=51, 0";  New = "376÷8=47, 0" },
    @{ Old = "584÷4=146, 0"; New = "774÷9=86, 0" },
    @{ Old = "555÷4=138, 3"; New = "249÷4=62, 1" },
    @{ Old = "733÷4=183, 1"; New = "138÷9=15, 3" },
    @{ Old = "996÷4=249, 0"; New = "262÷2=131, 0" },
    @{ Old = "253÷2=126, 1"; New = "557÷6=92, 5" },
    @{ Old = "371÷6=61, 5";  New = "844÷6=140, 4" },
    @{ Old = "938÷2=469, 0"; New = "811÷9=90, 1" },
    @{ Old = "122÷2=61, 0";  New = "777÷7=111, 0" },
    @{ Old = "423÷2=211, 1"; New = "801÷6=133, 3" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
